$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the cryptos.xlsx data refresh diff.
# Every target value in this sheet is stored as literal text (the source
# data keeps things like "56.852.15" / "0.0000156" / "1.00" as strings, not
# numbers), so numeric-looking values are pre-formatted as Text ("@") before
# being written -- this avoids Excel silently re-typing them as numbers and
# dropping significant trailing zeros / leading zeros / exact decimal text.

$ws.Range("D2").Value = "56.852.15"
$ws.Range("E2").Value = "  +4.03%  "
$ws.Range("D3").Value = "3.006.37"
$ws.Range("E3").Value = "  +3.97%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "511.20"
$ws.Range("E5").Value = "  +8.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.44"
$ws.Range("E6").Value = "  +11.02%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.434"
$ws.Range("E8").Value = "  +6.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.61"
$ws.Range("E9").Value = "  +14.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.110"
$ws.Range("E10").Value = "  +11.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.355"
$ws.Range("E11").Value = "  +6.56%  "
$ws.Range("E12").Value = "  +5.02%  "
$ws.Range("D13").Value = "3.514.42"
$ws.Range("E13").Value = "  +4.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.81"
$ws.Range("E14").Value = "  +10.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000156"
$ws.Range("E15").Value = "  +16.68%  "
$ws.Range("D16").Value = "56.845.00"
$ws.Range("E16").Value = "  +4.09%  "
$ws.Range("D17").Value = "3.001.36"
$ws.Range("E17").Value = "  +4.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.90"
$ws.Range("E18").Value = "  +8.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.49"
$ws.Range("E19").Value = "  +8.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.89"
$ws.Range("E20").Value = "  +11.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "330.54"
$ws.Range("E21").Value = "  +8.74%  "
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.485"
$ws.Range("E23").Value = "  +9.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.77"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.173"
$ws.Range("E25").Value = "  +14.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D27").Value = "0.0₃0918"
$ws.Range("E27").Value = "  +13.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.68"
$ws.Range("E28").Value = "  +8.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.26"
$ws.Range("E29").Value = "  +15.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.28"
$ws.Range("E30").Value = "  +15.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.79"
$ws.Range("E31").Value = "  +10.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.77"
$ws.Range("E32").Value = "  +9.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "156.42"
$ws.Range("E33").Value = "  +9.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.58"
$ws.Range("E34").Value = "  +8.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.68"
$ws.Range("E35").Value = "  +4.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.28"
$ws.Range("E36").Value = "  +5.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0683"
$ws.Range("E37").Value = "  +10.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.06"
$ws.Range("E38").Value = "  +6.30%  "
$ws.Range("D39").Value = "3.036.08"
$ws.Range("E39").Value = "  +4.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.96"
$ws.Range("E40").Value = "  +4.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.647"
$ws.Range("E42").Value = "  +7.09%  "
$ws.Range("D43").Value = "2.274.07"
$ws.Range("E43").Value = "  +11.67%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.42"
$ws.Range("E44").Value = "  +6.84%  "
$ws.Range("B45").Value = "ONDO"
$ws.Range("C45").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +4.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.65"
$ws.Range("E46").Value = "  +5.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.05"
$ws.Range("E47").Value = "  +27.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0240"
$ws.Range("E48").Value = "  +11.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.81"
$ws.Range("E49").Value = "  +8.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.41"
$ws.Range("E50").Value = "  +8.77%  "
$ws.Range("E51").Value = "  +9.33%  "
